$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 144
$ws.Range("I2").Value = 335
$ws.Range("J2").Value = 1438
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 417
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = 288
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 16
$ws.Range("S2").Value = 145
$ws.Range("T2").Value = 235
$ws.Range("U2").Value = 21
$ws.Range("V2").Value = 2254
$ws.Range("X2").Value = 2280
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 32
$ws.Range("AA2").Value = 13
